# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# de-de/zh-cn handback packages have been generated:
#  - Status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" everywhere it is used.
#  - The per-language "Latest Target File", "Latest Handback File" and
#    "Latest Handback DateTime" columns get populated for both data rows on
#    the zh-cn and de-de sheets, each with a hyperlink on the target file.
#  - A couple of columns are widened to fit the newly-written values.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc1f3ef073f772ba7c76e6e4f7f6b2872a4dce29/e2e/a.md"
$bMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc1f3ef073f772ba7c76e6e4f7f6b2872a4dce29/e2e/b.md"

# Cornflower blue (FF6495ED), same color already used by the existing
# hyperlink style in this workbook, expressed as an OLE/VB RGB value.
$hyperlinkColor = 15570276

# ----------------------------------------------------------------------
# Overview sheet: refresh the status text shown for each language column.
# ----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de columns to fit the longer status string.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ----------------------------------------------------------------------
# Helper that applies the handback info to a language sheet (zh-cn / de-de)
# ----------------------------------------------------------------------
function Set-HandbackInfo {
    param($ws, $handbackFileName, $handbackDateTime)

    # Status column
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Handback File / DateTime
    $ws.Range("J2").Value = $handbackFileName
    $ws.Range("J3").Value = $handbackFileName
    $ws.Range("K2").Value = $handbackDateTime
    $ws.Range("K3").Value = $handbackDateTime

    # Latest Target File (hyperlinked, same visual style as column A)
    $ws.Range("I2").Value = "a.md"
    $ws.Range("I3").Value = "a.md"

    # Rebuild the hyperlinks collection in the same order the report
    # generator writes it: A2, I2, A3, I3.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $aMdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $aMdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $bMdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "b.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $aMdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")

    # Hyperlinks.Add() re-stamps its own (theme-colored) style onto every
    # cell it touches, so re-apply the workbook's original custom
    # hyperlink look (underline + cornflower blue) afterwards, to each
    # cell individually so they all resolve to one consistent style.
    $ws.Range("A2").Font.Underline = 2
    $ws.Range("A2").Font.Color = $hyperlinkColor
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = $hyperlinkColor
    $ws.Range("A3").Font.Underline = 2
    $ws.Range("A3").Font.Color = $hyperlinkColor
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = $hyperlinkColor

    # Widen the Status and Latest Handback File columns to fit their
    # longer contents.
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

$zhcn = $wb.Worksheets.Item("zh-cn")
Set-HandbackInfo $zhcn "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-28 20:37:28"

$dede = $wb.Worksheets.Item("de-de")
Set-HandbackInfo $dede "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-28 20:37:35"
